# Automatic update of files.
# Increment the "Förändrad" (Changed) date in column C by one day
# for rows 2 through 5 on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 45244
$ws.Range("C3").Value = 45244
$ws.Range("C4").Value = 45244
$ws.Range("C5").Value = 45244
